# Updated via Streamlit Approval System
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ADJUSTMENT_AMOUNT (column AE) values for rows 2-19
$adjustmentAmounts = @{
    2  = 500
    3  = 1500
    4  = 7419
    5  = 100000
    6  = 1000
    7  = 30000
    8  = 6500
    9  = 3500
    10 = 5000
    11 = 400
    12 = 98894
    13 = 281859.52
    14 = 345
    15 = 277842
    16 = 100000
    17 = 1000
    18 = 500
    19 = 1500
}

foreach ($row in $adjustmentAmounts.Keys) {
    # AE<row>: ADJUSTMENT_AMOUNT, numeric
    $ws.Range("AE$row").Value = $adjustmentAmounts[$row]

    # AK<row>..AO<row>: COST_CENTER, LEDGER_NAME, LEDGER_UNDER, TO, BY
    # filled in with text "0" (stored as text, not a number)
    foreach ($col in @("AK", "AL", "AM", "AN", "AO")) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = "@"
        $cell.Value = "0"
    }
}
